# Auto-generated edit script: updates per commit 'chore: update Sheets via scheduled runner'
$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H62").Value = 18522518
$ws.Range("J62").Value = 2000
$ws.Range("L62").Value = 2000
$ws.Range("N62").Value = -3248
$ws.Range("H65").Value = 18522518
$ws.Range("J65").Value = 2000
$ws.Range("L65").Value = 10000
$ws.Range("N65").Value = -16240
$ws.Range("H132").Value = 10107443
$ws.Range("I132").Value = 17553210
$ws.Range("K132").Value = 52659630
$ws.Range("M132").Value = -52657100
$ws.Range("H137").Value = 1350.5319
$ws.Range("I137").Value = 976.2174
$ws.Range("J137").Value = 1709.25
$ws.Range("K137").Value = 2928.6522
$ws.Range("L137").Value = 5127.75
$ws.Range("M137").Value = -378.6522
$ws.Range("N137").Value = -10227.75

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H32").Value = 4112.731
$ws.Range("I32").Value = 4196.68
$ws.Range("J32").Value = 2014
$ws.Range("K32").Value = 4196.68
$ws.Range("L32").Value = 2014
$ws.Range("M32").Value = -3909.68
$ws.Range("N32").Value = -2588
$ws.Range("H46").Value = 4364
$ws.Range("I46").Value = 9991
$ws.Range("J46").Value = 3660.625
$ws.Range("K46").Value = 9991
$ws.Range("L46").Value = 3660.625
$ws.Range("M46").Value = -9672
$ws.Range("N46").Value = -4298.625
$ws.Range("H61").Value = 250001500
$ws.Range("I61").Value = 1000000000
$ws.Range("J61").Value = 2009.3334
$ws.Range("K61").Value = 1000000000
$ws.Range("L61").Value = 2009.3334
$ws.Range("M61").Value = -999999788
$ws.Range("N61").Value = -2433.3334
$ws.Range("H74").Value = 2862
$ws.Range("I74").Value = 2083.5
$ws.Range("K74").Value = 2083.5
$ws.Range("M74").Value = -1209.5
$ws.Range("H77").Value = 2862
$ws.Range("I77").Value = 2083.5
$ws.Range("K77").Value = 10417.5
$ws.Range("M77").Value = -6049.5
$ws.Range("H132").Value = 3031.2144
$ws.Range("J132").Value = 3753.25
$ws.Range("L132").Value = 11259.75
$ws.Range("N132").Value = -16319.75
$ws.Range("H136").Value = 250001500
$ws.Range("I136").Value = 1000000000
$ws.Range("J136").Value = 2009.3334
$ws.Range("K136").Value = 3000000000
$ws.Range("L136").Value = 6028.0002
$ws.Range("M136").Value = -2999997450
$ws.Range("N136").Value = -11128.0002

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H134").Value = 8643.333000000001
$ws.Range("I134").Value = 1275.25
$ws.Range("J134").Value = 17064
$ws.Range("K134").Value = 3825.75
$ws.Range("L134").Value = 51192
$ws.Range("M134").Value = -1290.75
$ws.Range("N134").Value = -56262

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H31").Value = 1333.591
$ws.Range("I31").Value = 1266.0555
$ws.Range("J31").Value = 1637.5
$ws.Range("K31").Value = 1266.0555
$ws.Range("L31").Value = 1637.5
$ws.Range("M31").Value = -971.0554999999999
$ws.Range("N31").Value = -2227.5
$ws.Range("H34").Value = 1333.591
$ws.Range("I34").Value = 1266.0555
$ws.Range("J34").Value = 1637.5
$ws.Range("K34").Value = 1266.0555
$ws.Range("L34").Value = 1637.5
$ws.Range("M34").Value = -1064.0555
$ws.Range("N34").Value = -2041.5

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H56").Value = 5612.8823
$ws.Range("I56").Value = 5612.8823
$ws.Range("K56").Value = 5612.8823
$ws.Range("M56").Value = -5082.8823
$ws.Range("H130").Value = 2410.2222
$ws.Range("J130").Value = 2582.75
$ws.Range("L130").Value = 7748.25
$ws.Range("N130").Value = -17788.25
$ws.Range("H131").Value = 18521932
$ws.Range("J131").Value = 3947.1304
$ws.Range("L131").Value = 11841.3912
$ws.Range("N131").Value = -21921.3912

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H64").Value = 15500
$ws.Range("J64").Value = 15500
$ws.Range("L64").Value = 15500
$ws.Range("N64").Value = -15996
$ws.Range("H67").Value = 15500
$ws.Range("J67").Value = 15500
$ws.Range("L67").Value = 15500
$ws.Range("N67").Value = -17216

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H26").Value = 2502.25
$ws.Range("J26").Value = 2500
$ws.Range("L26").Value = 2500
$ws.Range("N26").Value = -3090
$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("H47").Value = 7000
$ws.Range("J47").Value = 7000
$ws.Range("L47").Value = 7000
$ws.Range("N47").Value = -7980
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("H51").Value = 8000
$ws.Range("J51").Value = 8000
$ws.Range("L51").Value = 8000
$ws.Range("N51").Value = -8956
$ws.Range("H52").Value = 7000
$ws.Range("J52").Value = 7000
$ws.Range("L52").Value = 7000
$ws.Range("N52").Value = -7466
$ws.Range("H53").Value = 10051
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("H54").Value = 14000
$ws.Range("J54").Value = 14000
$ws.Range("L54").Value = 14000
$ws.Range("N54").Value = -15288
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("H58").Value = 5023.25
$ws.Range("I58").Value = 3046.5
$ws.Range("J58").Value = 7000
$ws.Range("K58").Value = 3046.5
$ws.Range("L58").Value = 7000
$ws.Range("M58").Value = -2786.5
$ws.Range("N58").Value = -7520
$ws.Range("H68").Value = 1715.5217
$ws.Range("I68").Value = 1693.85
$ws.Range("J68").Value = 1860
$ws.Range("K68").Value = 1693.85
$ws.Range("L68").Value = 1860
$ws.Range("M68").Value = -944.8499999999999
$ws.Range("N68").Value = -3358
$ws.Range("H71").Value = 1715.5217
$ws.Range("I71").Value = 1693.85
$ws.Range("J71").Value = 1860
$ws.Range("K71").Value = 8469.25
$ws.Range("L71").Value = 9300
$ws.Range("M71").Value = -4725.25
$ws.Range("N71").Value = -16788
$ws.Range("H123").Value = 40950
$ws.Range("J123").Value = 40950
$ws.Range("L123").Value = 40950
$ws.Range("N123").Value = -50750
$ws.Range("H132").Value = 2889.7827
$ws.Range("I132").Value = 3008.2
$ws.Range("K132").Value = 9024.599999999999
$ws.Range("M132").Value = -6494.599999999999
$ws.Range("H136").Value = 1687
$ws.Range("I136").Value = 802.25
$ws.Range("K136").Value = 2406.75
$ws.Range("M136").Value = 143.25
$ws.Range("N41").Value = ""
$ws.Range("M42").Value = ""
$ws.Range("M49").Value = ""
$ws.Range("M53").Value = ""
$ws.Range("N57").Value = ""

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("H132").Value = 5782.6113
$ws.Range("I132").Value = 8029.4
$ws.Range("K132").Value = 24088.2
$ws.Range("M132").Value = -21558.2
$ws.Range("M18").Value = ""
